# Generate Report for Handoff
# Update the localization-status report to reflect that b.md has been
# re-handed-off: new handoff xlf files were generated for zh-cn and de-de,
# and the overview / per-language sheets are updated accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: b.md row (row 3) moves from "Handed back: in sync with
# en-US" to "Ready for handoff", with an updated generation timestamp.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2017-02-09 15:09:04"

# ---------------------------------------------------------------------
# zh-cn sheet: b.md row (row 3) gets a freshly generated handoff file.
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("F3").Value = "False"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2017-02-09 15:08:44"
$zhcn.Range("R3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8491ec969280e038e6b8b5f820008e1b9b3d472f/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4aefb93b519c0ce6a27f4684a6760bab8c856e84/e2e/b.md."
$zhcn.Columns.Item(18).ColumnWidth = 39.17

# ---------------------------------------------------------------------
# de-de sheet: b.md row (row 3) gets a freshly generated handoff file.
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "False"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2017-02-09 15:09:04"
$dede.Range("R3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8491ec969280e038e6b8b5f820008e1b9b3d472f/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4aefb93b519c0ce6a27f4684a6760bab8c856e84/e2e/b.md."
$dede.Columns.Item(18).ColumnWidth = 39.17
